$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the RANK column (column A) content entirely - clear cells, don't shift others.
$ws.Range("A1:A4").ClearContents()

# Update MATCH_POINTS values (column B) to the new match-points based values.
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 0

$ws.Range("H9").Select()
